# Auto-generated edit script: update F (and some G) numeric counters
# per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 3615
$ws.Range("F8").Value = 2306
$ws.Range("F11").Value = 7673
$ws.Range("F12").Value = 7861
$ws.Range("F13").Value = 5036
$ws.Range("F16").Value = 630
$ws.Range("G16").Value = 70
$ws.Range("F17").Value = 5393
$ws.Range("G17").Value = 80
$ws.Range("F20").Value = 146
$ws.Range("F22").Value = 965
$ws.Range("F23").Value = 1505
$ws.Range("F24").Value = 2114
$ws.Range("F27").Value = 263
$ws.Range("F28").Value = 1104
$ws.Range("F30").Value = 760
$ws.Range("F31").Value = 57
$ws.Range("F32").Value = 800
$ws.Range("F33").Value = 1303
$ws.Range("F34").Value = 473
$ws.Range("F37").Value = 235

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 7812
$ws.Range("F9").Value = 21
$ws.Range("F28").Value = 81

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F7").Value = 695
$ws.Range("F9").Value = 9493
$ws.Range("F10").Value = 1820
$ws.Range("F12").Value = 122
$ws.Range("F15").Value = 307
$ws.Range("F16").Value = 2617
$ws.Range("F17").Value = 297
$ws.Range("F19").Value = 575

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 3615
$ws.Range("F6").Value = 695
$ws.Range("F7").Value = 1820
$ws.Range("F9").Value = 307
$ws.Range("F10").Value = 2617
$ws.Range("F11").Value = 297
$ws.Range("F12").Value = 5036
$ws.Range("F13").Value = 630
$ws.Range("G13").Value = 70
$ws.Range("F15").Value = 146
$ws.Range("F17").Value = 965
$ws.Range("F18").Value = 1505
$ws.Range("F20").Value = 575
$ws.Range("F21").Value = 575
$ws.Range("F22").Value = 21
$ws.Range("F26").Value = 263
$ws.Range("F28").Value = 760
$ws.Range("F29").Value = 57
$ws.Range("F30").Value = 800
$ws.Range("F32").Value = 1303
$ws.Range("F35").Value = 473
$ws.Range("F39").Value = 235
